# Adding the changes we made on may 9th
#
# 1) Insert 7 new accelerometer readings right after the header row
#    (these become the new rows 2-8, pushing the former rows 2-21 down to 9-28).
# 2) Append 3 more new accelerometer readings after the (now shifted) last
#    original row, becoming rows 29-31.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: insert 7 blank rows right below the header row and fill them ---
$ws.Rows("2:8").Insert()
# Inserting copies formatting from the row above (the bold header style);
# the source data rows carry no explicit style, so clear it back to default.
$ws.Range("A2:C8").ClearFormats()

$newTopRows = @(
    @(-2.288097732945488,  -4.339917659759537,  -3.936509232772027),
    @(-3.323771476745605,  -0.8808293342590332, -3.401906013488769),
    @(-5.066901056390067,   1.504574901179274,  -3.703831045251145),
    @(-4.698070149672659,  -0.4452685556913691, -8.655812690132544),
    @(-6.691070581737309,  -6.563599611583482,  -8.057312513652606),
    @(-4.573285604778092, -17.30385228207236,   12.64209365844724),
    @(-1.509485269847658,  -5.039251478094737,   6.722280000385433)
)

for ($i = 0; $i -lt $newTopRows.Length; $i++) {
    $r = 2 + $i
    $ws.Cells.Item($r, 1).Value = $newTopRows[$i][0]
    $ws.Cells.Item($r, 2).Value = $newTopRows[$i][1]
    $ws.Cells.Item($r, 3).Value = $newTopRows[$i][2]
}

# --- Step 2: append 3 more rows after the existing data (new rows 29-31) ---
$newBottomRows = @(
    @(-3.77490947121079,  -6.046053083319467, -15.12075865896126),
    @(3.047089124980793,  -31.00068785014912,   2.002202786897417),
    @(7.431886748263759,  -15.14435035304032,   4.078887035972246)
)

for ($i = 0; $i -lt $newBottomRows.Length; $i++) {
    $r = 29 + $i
    $ws.Cells.Item($r, 1).Value = $newBottomRows[$i][0]
    $ws.Cells.Item($r, 2).Value = $newBottomRows[$i][1]
    $ws.Cells.Item($r, 3).Value = $newBottomRows[$i][2]
}
